$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace numeric category codes in column A with descriptive text labels
# 1 -> "TOYS", 2 -> "AGD"
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $curVal = $cell.Value2
    if ($curVal -eq 1) {
        $cell.Value = "TOYS"
    } elseif ($curVal -eq 2) {
        $cell.Value = "AGD"
    }
}

# Reset zoom to 100% (removes custom zoom scaling) and move the selection to A2
$excel.ActiveWindow.Zoom = 100
$ws.Range("A2").Select()

# Give column A an explicit custom width
$ws.Columns.Item(1).ColumnWidth = 8

$wb.Save()
